$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank columns before column Q (17).
# This shifts the old "Comments" column (Q) to T, and the old extra column (R) to U,
# while preserving all existing data/styles in those columns.
$ws.Range("Q1:S1").EntireColumn.Insert()

# Rename the old "special1"/"special2" headers and label the newly inserted columns.
$ws.Range("O1").Value = "stars"
$ws.Range("P1").Value = "filter"
$ws.Range("Q1").Value = "variants"
$ws.Range("R1").Value = "unused1"
$ws.Range("S1").Value = "unused2"

# Tag the rows belonging to the tarsp2005 variant (rows with Fase = 7) in the new
# "variants" column.
$tarsp2005Rows = @(4, 56, 60, 74, 75, 84, 98, 104, 110, 139, 153)
foreach ($r in $tarsp2005Rows) {
    $ws.Cells.Item($r, 17).Value = "tarsp2005"
}

# Re-apply the AutoFilter over the new range, filtering on column J (Fase = 7),
# which hides every row that isn't part of the tarsp2005 variant.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:U155").AutoFilter(10, @("7"), 7)

# Keep the hidden _FilterDatabase defined name in sync with the new filter range.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$U`$155"
